# Horarios actualizados Línea 141 - 1007
$wb = $excel.ActiveWorkbook

$newTime = "03:00:58"

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newTime"
$ws1.Range("A3").Value = "Total filas: 5"

# Update existing data rows 6-9: refresh scrape time + minutes
$ws1.Range("A6").Value = $newTime
$ws1.Range("D6").Value = 1

$ws1.Range("A7").Value = $newTime
$ws1.Range("B7").Value = "03:47"
$ws1.Range("D7").Value = 47

$ws1.Range("A8").Value = $newTime
$ws1.Range("D8").Value = 61

$ws1.Range("A9").Value = $newTime
$ws1.Range("D9").Value = 107

# New row 10
$ws1.Range("A10").Value = $newTime
$ws1.Range("B10").Value = "04:52"
$ws1.Range("C10").Value = "11_ETCHEVERRY"
$ws1.Range("D10").Value = 112
$ws1.Range("E10").Value = "LP1912"

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: $newTime"

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: $newTime"
